$d = $word.ActiveDocument

$d.Content.Find.Execute("430÷4=107, 2", $false, $false, $false, $false, $false, $true, 1, $false, "432÷2=216, 0", 2) | Out-Null
$d.Content.Find.Execute("389÷7=55, 4", $false, $false, $false, $false, $false, $true, 1, $false, "455÷7=65, 0", 2) | Out-Null
$d.Content.Find.Execute("656÷6=109, 2", $false, $false, $false, $false, $false, $true, 1, $false, "735÷9=81, 6", 2) | Out-Null
$d.Content.Find.Execute("746÷6=124, 2", $false, $false, $false, $false, $false, $true, 1, $false, "494÷5=98, 4", 2) | Out-Null
$d.Content.Find.Execute("951÷2=475, 1", $false, $false, $false, $false, $false, $true, 1, $false, "540÷7=77, 1", 2) | Out-Null
$d.Content.Find.Execute("443÷8=55, 3", $false, $false, $false, $false, $false, $true, 1, $false, "671÷3=223, 2", 2) | Out-Null
$d.Content.Find.Execute("459÷8=57, 3", $false, $false, $false, $false, $false, $true, 1, $false, "752÷5=150, 2", 2) | Out-Null
$d.Content.Find.Execute("479÷3=159, 2", $false, $false, $false, $false, $false, $true, 1, $false, "957÷2=478, 1", 2) | Out-Null
$d.Content.Find.Execute("713÷9=79, 2", $false, $false, $false, $false, $false, $true, 1, $false, "740÷7=105, 5", 2) | Out-Null
$d.Content.Find.Execute("603÷7=86, 1", $false, $false, $false, $false, $false, $true, 1, $false, "781÷4=195, 1", 2) | Out-Null
$d.Content.Find.Execute("822÷3=274, 0", $false, $false, $false, $false, $false, $true, 1, $false, "795÷4=198, 3", 2) | Out-Null
$d.Content.Find.Execute("153÷9=17, 0", $false, $false, $false, $false, $false, $true, 1, $false, "587÷4=146, 3", 2) | Out-Null
$d.Content.Find.Execute("272÷7=38, 6", $false, $false, $false, $false, $false, $true, 1, $false, "136÷8=17, 0", 2) | Out-Null
$d.Content.Find.Execute("530÷9=58, 8", $false, $false, $false, $false, $false, $true, 1, $false, "136÷8=17, 0", 2) | Out-Null
$d.Content.Find.Execute("914÷3=304, 2", $false, $false, $false, $false, $false, $true, 1, $false, "104÷5=20, 4", 2) | Out-Null
$d.Content.Find.Execute("877÷7=125, 2", $false, $false, $false, $false, $false, $true, 1, $false, "517÷9=57, 4", 2) | Out-Null
$d.Content.Find.Execute("648÷8=81, 0", $false, $false, $false, $false, $false, $true, 1, $false, "268÷8=33, 4", 2) | Out-Null
$d.Content.Find.Execute("866÷8=108, 2", $false, $false, $false, $false, $false, $true, 1, $false, "530÷9=58, 8", 2) | Out-Null
$d.Content.Find.Execute("219÷2=109, 1", $false, $false, $false, $false, $false, $true, 1, $false, "918÷6=153, 0", 2) | Out-Null
$d.Content.Find.Execute("426÷9=47, 3", $false, $false, $false, $false, $false, $true, 1, $false, "240÷2=120, 0", 2) | Out-Null
$d.Content.Find.Execute("895÷3=298, 1", $false, $false, $false, $false, $false, $true, 1, $false, "150÷3=50, 0", 2) | Out-Null
$d.Content.Find.Execute("738÷5=147, 3", $false, $false, $false, $false, $false, $true, 1, $false, "327÷5=65, 2", 2) | Out-Null
$d.Content.Find.Execute("695÷7=99, 2", $false, $false, $false, $false, $false, $true, 1, $false, "242÷7=34, 4", 2) | Out-Null
$d.Content.Find.Execute("252÷9=28, 0", $false, $false, $false, $false, $false, $true, 1, $false, "581÷6=96, 5", 2) | Out-Null
$d.Content.Find.Execute("726÷2=363, 0", $false, $false, $false, $false, $false, $true, 1, $false, "160÷4=40, 0", 2) | Out-Null
